$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$arr = New-Object "object[,]" 24,1
$arr[0,0] = 4.836692403563232
$arr[1,0] = 4.757803138350983
$arr[2,0] = 4.709707529428421
$arr[3,0] = 4.690221876971522
$arr[4,0] = 4.686993973645599
$arr[5,0] = 4.709444242241663
$arr[6,0] = 4.809436698797802
$arr[7,0] = 5.007046127501917
$arr[8,0] = 5.151594405828754
$arr[9,0] = 5.24202854214806
$arr[10,0] = 5.277392942209253
$arr[11,0] = 5.269803666181264
$arr[12,0] = 5.244950173512259
$arr[13,0] = 5.229647656010439
$arr[14,0] = 5.147315346957807
$arr[15,0] = 5.10976091646931
$arr[16,0] = 5.088119540427556
$arr[17,0] = 5.080785881392673
$arr[18,0] = 5.113763086694421
$arr[19,0] = 5.252266745148467
$arr[20,0] = 5.35406034236956
$arr[21,0] = 5.300058621474275
$arr[22,0] = 5.111953861934731
$arr[23,0] = 4.953588568731957
$ws.Range("B2:B25").Value = $arr

$arr = New-Object "object[,]" 24,1
$arr[0,0] = 5.587913769660256
$arr[1,0] = 5.58204793906518
$arr[2,0] = 5.579191649810346
$arr[3,0] = 5.578216303340128
$arr[4,0] = 5.578065772815987
$arr[5,0] = 5.579177730729765
$arr[6,0] = 5.585737268980873
$arr[7,0] = 5.604459607570914
$arr[8,0] = 5.621709172680452
$arr[9,0] = 5.630297034069592
$arr[10,0] = 5.633653896176374
$arr[11,0] = 5.632926301974512
$arr[12,0] = 5.630571114417807
$arr[13,0] = 5.629142092957714
$arr[14,0] = 5.621162700097431
$arr[15,0] = 5.616456175354652
$arr[16,0] = 5.613818956416908
$arr[17,0] = 5.61293808731543
$arr[18,0] = 5.616949976004003
$arr[19,0] = 5.631260059274084
$arr[20,0] = 5.641222675289925
$arr[21,0] = 5.635850199798385
$arr[22,0] = 5.616726514849727
$arr[23,0] = 5.598775141453702
$ws.Range("D2:D25").Value = $arr

$arr = New-Object "object[,]" 24,1
$arr[0,0] = 12.20519040871745
$arr[1,0] = 12.17880831628767
$arr[2,0] = 12.1652800854924
$arr[3,0] = 12.16044157835812
$arr[4,0] = 12.15967894868909
$arr[5,0] = 12.16521209793162
$arr[6,0] = 12.19554148623213
$arr[7,0] = 12.27605639019627
$arr[8,0] = 12.34779977638117
$arr[9,0] = 12.38310873994961
$arr[10,0] = 12.39685739552211
$arr[11,0] = 12.39387967607317
$arr[12,0] = 12.38423231461405
$arr[13,0] = 12.37837204425603
$arr[14,0] = 12.34554549707147
$arr[15,0] = 12.32608760222873
$arr[16,0] = 12.31514790563803
$arr[17,0] = 12.31148737413398
$arr[18,0] = 12.32813289498893
$arr[19,0] = 12.38705577557071
$arr[20,0] = 12.42776473589656
$arr[21,0] = 12.40583859714636
$arr[22,0] = 12.32720744859571
$arr[23,0] = 12.25204468814137
$ws.Range("E2:E25").Value = $arr

$arr = New-Object "object[,]" 24,1
$arr[0,0] = 50.45908136372888
$arr[1,0] = 50.3194270853546
$arr[2,0] = 50.24705989510221
$arr[3,0] = 50.22094349030052
$arr[4,0] = 50.21681083906021
$arr[5,0] = 50.2466940109956
$arr[6,0] = 50.40815787479795
$arr[7,0] = 50.83039650524249
$arr[8,0] = 51.20400110284908
$arr[9,0] = 51.38743713527315
$arr[10,0] = 51.45880779298106
$arr[11,0] = 51.44335254520664
$arr[12,0] = 51.39327080573432
$arr[13,0] = 51.36284171112914
$arr[14,0] = 51.19228181949933
$arr[15,0] = 51.09108067407364
$arr[16,0] = 51.03414348128775
$arr[17,0] = 51.01508472017595
$arr[18,0] = 51.10172237712146
$arr[19,0] = 51.40792952557649
$arr[20,0] = 51.61915461321678
$arr[21,0] = 51.50541524794869
$arr[22,0] = 51.09690738390542
$arr[23,0] = 50.70495136807557
$ws.Range("F2:F25").Value = $arr

$arr = New-Object "object[,]" 24,1
$arr[0,0] = 3.746203674247984
$arr[1,0] = 3.750485996051044
$arr[2,0] = 3.753249572876031
$arr[3,0] = 3.754409639891246
$arr[4,0] = 3.754604318872058
$arr[5,0] = 3.753265080558034
$arr[6,0] = 3.747652445883862
$arr[7,0] = 3.737704722336679
$arr[8,0] = 3.731032691341991
$arr[9,0] = 3.728133747806806
$arr[10,0] = 3.727055432455359
$arr[11,0] = 3.727286803960819
$arr[12,0] = 3.728044645076201
$arr[13,0] = 3.728511373868751
$arr[14,0] = 3.731224873347915
$arr[15,0] = 3.732924305734387
$arr[16,0] = 3.733914601738834
$arr[17,0] = 3.734252106139328
$arr[18,0] = 3.732742071630501
$arr[19,0] = 3.727821521919185
$arr[20,0] = 3.724718975664513
$arr[21,0] = 3.726364538405882
$arr[22,0] = 3.732824418326373
$arr[23,0] = 3.740283422921848
$ws.Range("G2:G25").Value = $arr

$arr = New-Object "object[,]" 24,1
$arr[0,0] = 23.13577021670493
$arr[1,0] = 23.20308603613041
$arr[2,0] = 23.24782119517909
$arr[3,0] = 23.26690565098806
$arr[4,0] = 23.27012618634865
$arr[5,0] = 23.24807511658234
$arr[6,0] = 23.15827370735938
$arr[7,0] = 23.00922811053761
$arr[8,0] = 22.91629167557812
$arr[9,0] = 22.87762622528928
$arr[10,0] = 22.86350542243466
$arr[11,0] = 22.86652339445157
$arr[12,0] = 22.87645404828488
$arr[13,0] = 22.88260475285347
$arr[14,0] = 22.91889136806201
$arr[15,0] = 22.94207822761199
$arr[16,0] = 22.95575455972497
$arr[17,0] = 22.96044345438761
$arr[18,0] = 22.93957476074959
$arr[19,0] = 22.87352302295849
$arr[20,0] = 22.83339171921053
$arr[21,0] = 22.85453210628723
$arr[22,0] = 22.94070550084435
$arr[23,0] = 23.0466449080885
$ws.Range("I2:I25").Value = $arr

$arr = New-Object "object[,]" 24,1
$arr[0,0] = 10.57228521820152
$arr[1,0] = 10.5828884959858
$arr[2,0] = 10.59129938261032
$arr[3,0] = 10.59520406605322
$arr[4,0] = 10.59588123348329
$arr[5,0] = 10.59135011147883
$arr[6,0] = 10.57554640408179
$arr[7,0] = 10.55966234293569
$arr[8,0] = 10.55723503440674
$arr[9,0] = 10.55814238175526
$arr[10,0] = 10.55877534059786
$arr[11,0] = 10.55862615168593
$arr[12,0] = 10.55818865619962
$arr[13,0] = 10.55795836246554
$arr[14,0] = 10.5572162269355
$arr[15,0] = 10.55727627602919
$arr[16,0] = 10.55750015279238
$arr[17,0] = 10.55760846703575
$arr[18,0] = 10.55725028670479
$arr[19,0] = 10.55830930573277
$arr[20,0] = 10.56068803582641
$arr[21,0] = 10.55926413850581
$arr[22,0] = 10.55726144669341
$arr[23,0] = 10.56233817608212
$ws.Range("J2:J25").Value = $arr

$arr = New-Object "object[,]" 24,1
$arr[0,0] = 22.34956385378855
$arr[1,0] = 22.06844570978598
$arr[2,0] = 21.90069366044535
$arr[3,0] = 21.83363330331894
$arr[4,0] = 21.82257874709569
$arr[5,0] = 21.89978389545979
$arr[6,0] = 22.2516747604947
$arr[7,0] = 22.97666600480498
$arr[8,0] = 23.52568472387361
$arr[9,0] = 23.77794025364134
$arr[10,0] = 23.87373463266882
$arr[11,0] = 23.85309296999602
$arr[12,0] = 23.78581639746236
$arr[13,0] = 23.74464017344452
$arr[14,0] = 23.50924195165999
$arr[15,0] = 23.36541017216481
$arr[16,0] = 23.2829229363455
$arr[17,0] = 23.25503831064338
$arr[18,0] = 23.38069701964309
$arr[19,0] = 23.80557051744237
$arr[20,0] = 24.08478722982339
$arr[21,0] = 23.93565252672449
$arr[22,0] = 23.37378519502697
$arr[23,0] = 22.77730171147421
$ws.Range("K2:K25").Value = $arr

$arr = New-Object "object[,]" 24,1
$arr[0,0] = 21.12225618396107
$arr[1,0] = 21.18762663375729
$arr[2,0] = 21.22972809481873
$arr[3,0] = 21.2473795333621
$arr[4,0] = 21.2503404512585
$arr[5,0] = 21.22996414387608
$arr[6,0] = 21.14438891552619
$arr[7,0] = 20.99211313806057
$arr[8,0] = 20.88964445096642
$arr[9,0] = 20.84505846361287
$arr[10,0] = 20.82846555733974
$arr[11,0] = 20.83202621254532
$arr[12,0] = 20.84368753039204
$arr[13,0] = 20.85086827931787
$arr[14,0] = 20.89259898732713
$arr[15,0] = 20.91871820613306
$arr[16,0] = 20.93393217224859
$arr[17,0] = 20.93911616620279
$arr[18,0] = 20.91591802039695
$arr[19,0] = 20.84025443094714
$arr[20,0] = 20.79249882209527
$arr[21,0] = 20.81783203406028
$arr[22,0] = 20.91718336858687
$arr[23,0] = 21.03165107535158
$ws.Range("N2:N25").Value = $arr

